# August 5 forecast update
# The NTC (Named Tropical Cyclone) adjustment value on the County sheet
# (cell B1) changes from 170 to 143. All dependent formulas across the
# County, State, Province, and Caribbean-Central America sheets reference
# County!$B$1 and will recalculate automatically.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("County")
$ws.Range("B1").Value = 143
